# Update receptor/edge expression values (columns M:T) for rows 2-5 with
# newly computed TPM-based figures, per "update scripts wuth new tpm".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 5.474275674611111
$ws.Range("R2").Value = 49.26848107149999
$ws.Range("S2").Value = 0.2176904746803693
$ws.Range("T2").Value = 0.2176904746803693

# Row 3
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83271999999999
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 9.279565548888888
$ws.Range("R3").Value = 83.51608993999999
$ws.Range("S3").Value = 0.3690119294748028
$ws.Range("T3").Value = 0.3690119294748029

# Row 4
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 8.827674060416665
$ws.Range("R4").Value = 79.44906654374999
$ws.Range("S4").Value = 0.3510419771967738
$ws.Range("T4").Value = 0.3510419771967739

# Row 5
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 1.565545847944444
$ws.Range("R5").Value = 14.0899126315
$ws.Range("S5").Value = 0.06225561864805391
$ws.Range("T5").Value = 0.06225561864805392
